# Lecture partielle de l'EDT M1 MIAGE.
# Update weekday labels and their corresponding dates to reflect a later week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekday name cells (shared strings): jeudi -> lundi, vendredi -> mardi, samedi -> mercredi
$ws.Range("B2").Value = "lundi"
$ws.Range("B5").Value = "mardi"
$ws.Range("B8").Value = "mercredi"

# Corresponding date cells, keep existing date formatting, change the underlying date value
$ws.Range("A2").Value = (Get-Date -Year 2026 -Month 1 -Day 5 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("A5").Value = (Get-Date -Year 2026 -Month 1 -Day 13 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("A8").Value = (Get-Date -Year 2026 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0).Date
